$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws 'D2' '57.944.60'
$ws.Range('E2').Value = '  +2.78%  '
Set-TextValue $ws 'D3' '2.325.92'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.36%  '
Set-TextValue $ws 'D5' '541.75'
$ws.Range('E5').Value = '  +5.92%  '
Set-TextValue $ws 'D6' '134.61'
$ws.Range('E6').Value = '  +1.91%  '
Set-TextValue $ws 'D7' '1.00'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +0.56%  '
Set-TextValue $ws 'D9' '2.354.69'
$ws.Range('E9').Value = '  +0.87%  '
Set-TextValue $ws 'D10' '0.102'
$ws.Range('E10').Value = '  +2.22%  '
$ws.Range('E11').Value = '  +1.09%  '
Set-TextValue $ws 'D12' '5.37'
$ws.Range('E12').Value = '  +1.49%  '
Set-TextValue $ws 'D13' '0.354'
$ws.Range('E13').Value = '  +4.71%  '
Set-TextValue $ws 'D14' '2.775.76'
$ws.Range('E14').Value = '  +1.24%  '
Set-TextValue $ws 'D15' '23.48'
$ws.Range('E15').Value = '  -0.34%  '
Set-TextValue $ws 'D16' '57.865.73'
$ws.Range('E16').Value = '  +2.67%  '
Set-TextValue $ws 'D17' '0.0000133'
$ws.Range('E17').Value = '  +0.81%  '
Set-TextValue $ws 'D18' '2.352.80'
$ws.Range('E18').Value = '  +0.61%  '
Set-TextValue $ws 'D19' '337.19'
$ws.Range('E19').Value = '  +4.31%  '
Set-TextValue $ws 'D20' '10.53'
$ws.Range('E20').Value = '  +1.08%  '
Set-TextValue $ws 'D21' '4.20'
$ws.Range('E21').Value = '  +1.62%  '
Set-TextValue $ws 'D22' '6.77'
$ws.Range('E22').Value = '  +2.08%  '
Set-TextValue $ws 'D23' '0.997'
$ws.Range('E23').Value = '  +0.02%  '
Set-TextValue $ws 'D24' '62.16'
$ws.Range('E24').Value = '  +1.01%  '
Set-TextValue $ws 'D25' '0.169'
$ws.Range('E25').Value = '  +4.34%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws 'D26' '0.998'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D27' '8.45'
$ws.Range('E27').Value = '  -2.75%  '
Set-TextValue $ws 'D28' '1.42'
$ws.Range('E28').Value = '  +9.00%  '
$ws.Range('E29').Value = '  +4.94%  '
Set-TextValue $ws 'D30' '170.57'
$ws.Range('E30').Value = '  +1.89%  '
Set-TextValue $ws 'D31' '0.0₃0735'
$ws.Range('E31').Value = '  +2.18%  '
Set-TextValue $ws 'D32' '6.15'
$ws.Range('E32').Value = '  +0.92%  '
Set-TextValue $ws 'D33' '18.55'
$ws.Range('E33').Value = '  +1.61%  '
Set-TextValue $ws 'D34' '1.01'
$ws.Range('E34').Value = '  +14.14%  '
Set-TextValue $ws 'D35' '0.999'
$ws.Range('E35').Value = '  -0.01%  '
Set-TextValue $ws 'D36' '0.999'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D37' '1.25'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D38' '4.11'
$ws.Range('E38').Value = '  +4.60%  '
$ws.Range('E39').Value = '  +4.01%  '
Set-TextValue $ws 'D40' '39.34'
$ws.Range('E40').Value = '  +2.31%  '
Set-TextValue $ws 'D41' '148.46'
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  +0.18%  '
Set-TextValue $ws 'D43' '3.62'
$ws.Range('E43').Value = '  +1.48%  '
Set-TextValue $ws 'D44' '282.09'
$ws.Range('E44').Value = '  +1.25%  '
Set-TextValue $ws 'D45' '0.0930'
$ws.Range('E45').Value = '  +0.60%  '
Set-TextValue $ws 'D46' '19.19'
$ws.Range('E46').Value = '  +7.15%  '
Set-TextValue $ws 'D47' '0.0506'
$ws.Range('E47').Value = '  +2.07%  '
Set-TextValue $ws 'D48' '0.559'
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('B50').Value = 'Polygon'
$ws.Range('C50').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D50' '0.383'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D51' '17.48'
$ws.Range('E51').Value = '  +1.96%  '
